$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D8")
$tmp  = $ws.Range("AZ9999")

# D8 currently carries a "quote-prefixed text" cell style (s=11). Writing a new
# formula into the cell directly would make Excel drop the quote-prefix flag and
# re-derive a (new) cell style. To keep D8 on its original style, stash a copy of
# its current formatting on a scratch cell first, apply the new formula to D8, and
# then restore the original formatting from the scratch cell.
$cell.Copy()
$tmp.PasteSpecial(-4122)  # xlPasteFormats

$cell.Formula = '=SUBSTITUTE(LEFT(CELL("filename",A1),FIND("[",CELL("filename",A1),1)-1),"\XLS\","\XML\")'

$tmp.Copy()
$cell.PasteSpecial(-4122)  # xlPasteFormats

$tmp.Clear()
